$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(16.51415364256055,9.866888820131091,0,15.61786258095976,36.45562086513911,3.638166415003166,0,23.03879365177308,7.492013016042842,0,12.99190345555264,0,17.538709072854,22.89726700555503),
    @(16.04284133325686,9.68681919019834,0,15.64192260551311,36.45947151408759,3.640186424805744,0,23.13801665184768,7.49225289297685,0,12.96375519622687,0,17.58382268671548,22.95399554411834),
    @(15.7483014481524,9.573450524237179,0,15.65867064770506,36.47088384310157,3.641493045685253,0,23.20377204819147,7.492554825750665,0,12.94823790946514,0,17.6133047162144,22.99445556302668),
    @(15.62715704738986,9.526582102406529,0,15.66599271063312,36.47780975328335,3.642042235250487,0,23.23178109125459,7.492716935820981,0,12.94236288382812,0,17.62576781088191,23.01235468317049),
    @(15.60697909178788,9.518760254658893,0,15.66723856949501,36.47909719377536,3.642134439813991,0,23.23650519035397,7.49274621913306,0,12.94141454332281,0,17.62786443644745,23.01541193672555),
    @(15.74667193069574,9.572821104011393,0,15.65876738230043,36.47096803667094,3.641500384435348,0,23.20414487828838,7.492556853592555,0,12.94815685576014,0,17.61347097909761,22.99469124791321),
    @(16.35281800323211,9.805403912711016,0,15.62574882714441,36.45507110607056,3.638849179966228,0,23.07200171499499,7.492063737943409,0,12.98183413465692,0,17.55389485750706,22.91565651076369),
    @(17.49334888618464,10.23775816491743,0,15.57665054265263,36.49562976478927,3.634174050702676,0,22.85129381913128,7.492315964944621,0,13.06168308740564,0,17.45117034936183,22.8054948212595),
    @(18.29314344867786,10.5390732704641,0,15.55009056583669,36.5689964915099,3.631055245434909,0,22.71267961271665,7.493234304374756,0,13.1284760189056,0,17.38424775923303,22.75209332098017),
    @(18.64709254617911,10.67224282206007,0,15.54006631842809,36.61176904498685,3.629704325584104,0,22.65475708019883,7.493809063753989,0,13.16055892101622,0,17.35564919860716,22.73381535583992),
    @(18.77958499996817,10.72208209889695,0,15.53656565983487,36.6293093126443,3.629202468656877,0,22.63356375453755,7.494049094864945,0,13.17294588725036,0,17.34508422596943,22.72776121161712),
    @(18.75112066541559,10.71137495225319,0,15.53730646614723,36.62547210712992,3.62931012156393,0,22.6380951352926,7.493996407530804,0,13.17026765750198,0,17.34734781952555,22.72902647134512),
    @(18.65802429405208,10.67635509222517,0,15.53977240258469,36.61318522496284,3.629662843243057,0,22.65299863412051,7.493828363805526,0,13.16157327440543,0,17.35477471074251,22.73329988290649),
    @(18.60079613101226,10.65482688449609,0,15.54132129771235,36.60583381598278,3.629880158181805,0,22.66222398888688,7.493728341192902,0,13.15627849743524,0,17.35935834774167,22.7360304868162),
    @(18.26980272093626,10.53028938424071,0,15.55078703359552,36.56638964754242,3.631144892210602,0,22.71656850014661,7.493199885194292,0,13.12641296497064,0,17.38615381338412,22.75340907612594),
    @(18.06412933433905,10.4528700108155,0,15.55712062979472,36.5445935749839,3.631938106647588,0,22.75122371777677,7.492915768729429,0,13.10852200058221,0,17.40306403841913,22.76561263353738),
    @(17.94490772098021,10.40797490512059,0,15.5609572986598,36.53294218514274,3.632400730922673,0,22.77163954909231,7.492767141970807,0,13.09839185730951,0,17.41296403466259,22.77319780528754),
    @(17.90438651567322,10.3927123058894,0,15.56228962432434,36.52914947114613,3.632558466228902,0,22.77863489152176,7.492719365519525,0,13.09498968032345,0,17.41634585219148,22.77586314573237),
    @(18.08612012193037,10.46114947703935,0,15.55642635950719,36.54682225023755,3.631853006801531,0,22.74748459737148,7.492944484153773,0,13.11040998260455,0,17.40124594510898,22.76425494533381),
    @(18.68541165420111,10.68665747959901,0,15.53904008807175,36.61675779748597,3.629558977270037,0,22.64860099379925,7.493877116433039,0,13.16412062191792,0,17.35258607411734,22.73202112219638),
    @(19.06805754316905,10.83059418641484,0,15.52939809335288,36.67028926830688,3.628116254657888,0,22.58829324015215,7.49461701470561,0,13.20060696185242,0,17.32232653499369,22.71601070619065),
    @(18.86469453274717,10.75409662507764,0,15.53438696495085,36.64100564194698,3.628881103622076,0,22.62008467661317,7.494210253000539,0,13.1810091136971,0,17.33833567708932,22.7240924292566),
    @(18.07618112153478,10.45740752953058,0,15.55673963041415,36.54581192608401,3.631891459895702,0,22.74917352105127,7.492931456053562,0,13.10955594074614,0,17.40206735032041,22.76486698336172),
    @(17.190916062557,10.12353550907505,0,15.58826015892516,36.4769890592993,3.635383060981254,0,22.9068757186085,7.492118256439911,0,13.03863229119492,0,17.47745526605713,22.8304747346705)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
